$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 165, pushing old rows 165..255 down to 166..256.
$ws.Rows.Item(165).Insert()

# Populate the new row 165 with the new data point.
$ws.Range("A165").Value = 3
$ws.Range("B165").Value = 'Femacal de La Calera'
$ws.Range("C165").Value = 'Coquimbo'
$ws.Range("D165").Value = 45089
$ws.Range("E165").Value = 5
$ws.Range("F165").Value = 100112026
$ws.Range("G165").Value = 'Haba'
$ws.Range("H165").Value = 'Sin especificar'
$ws.Range("I165").Value = 'Primera'
$ws.Range("J165").Value = 38
$ws.Range("K165").Value = 18000
$ws.Range("L165").Value = 18000
$ws.Range("M165").Value = 18000
$ws.Range("N165").Value = '$/saco 25 kilos'
$ws.Range("O165").Value = 'Provincia de Quillota'
$ws.Range("P165").Value = 720
$ws.Range("Q165").Value = 25
$ws.Range("R165").Value = 'Hortaliza'
